$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 187, shifting existing rows 187:240 down to 188:241
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new weekly record
$ws.Range("A187").Value = 6
$ws.Range("B187").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C187").Value = "Metropolitana"
$ws.Range("D187").Value = 44809
$ws.Range("E187").Value = 13
$ws.Range("F187").Value = 100112001
$ws.Range("G187").Value = "Berenjena"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 298
$ws.Range("K187").Value = 9000
$ws.Range("L187").Value = 11000
$ws.Range("M187").Value = 9671
$ws.Range("N187").Value = "$/caja 50 unidades"
$ws.Range("O187").Value = "Región de Arica y Parinacota"
$ws.Range("P187").Value = 193
$ws.Range("Q187").Value = 50
$ws.Range("R187").Value = "Hortaliza"
